$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.534.61"
$ws.Range("E2").Value = "  +2.29%  "

$ws.Range("D3").Value = "3.629.87"
$ws.Range("E3").Value = "  +4.84%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "237.07"
$ws.Range("E5").Value = "  +0.84%  "

$ws.Range("D6").Value = "658.68"
$ws.Range("E6").Value = "  +5.34%  "

$ws.Range("D7").Value = "1.46"
$ws.Range("E7").Value = "  +2.34%  "

$ws.Range("D8").Value = "0.402"
$ws.Range("E8").Value = "  +3.15%  "

$ws.Range("E9").Value = "  -0.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.56%  "

$ws.Range("D11").Value = "3.627.93"
$ws.Range("E11").Value = "  +4.85%  "

$ws.Range("D12").Value = "0.201"

$ws.Range("D13").Value = "42.33"
$ws.Range("E13").Value = "  -2.72%  "

$ws.Range("D14").Value = "6.45"
$ws.Range("E14").Value = "  +3.47%  "

$ws.Range("D15").Value = "4.314.73"
$ws.Range("E15").Value = "  +5.18%  "

$ws.Range("D16").Value = "95.326.49"
$ws.Range("E16").Value = "  +2.21%  "

$ws.Range("D17").Value = "0.0000253"
$ws.Range("E17").Value = "  +2.01%  "

$ws.Range("D18").Value = "3.631.22"
$ws.Range("E18").Value = "  +5.00%  "

$ws.Range("D19").Value = "7.94"
$ws.Range("E19").Value = "  -4.15%  "

$ws.Range("D20").Value = "12.93"
$ws.Range("E20").Value = "  +9.54%  "

$ws.Range("D21").Value = "18.05"
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").Value = "3.55"
$ws.Range("E22").Value = "  +4.59%  "

$ws.Range("D23").Value = "0.484"
$ws.Range("E23").Value = "  -2.89%  "

$ws.Range("D24").Value = "504.74"
$ws.Range("E24").Value = "  +0.32%  "

$ws.Range("D25").Value = "0.0000197"
$ws.Range("E25").Value = "  +7.47%  "

$ws.Range("D26").Value = "6.64"
$ws.Range("E26").Value = "  -2.37%  "

$ws.Range("D27").Value = "91.82"
$ws.Range("E27").Value = "  -3.22%  "

$ws.Range("D28").Value = "3.817.66"
$ws.Range("E28").Value = "  +4.78%  "

$ws.Range("D29").Value = "12.55"
$ws.Range("E29").Value = "  +2.79%  "

$ws.Range("D30").Value = "3.11"

$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("D32").Value = "11.28"
$ws.Range("E32").Value = "  -0.43%  "

$ws.Range("D33").Value = "0.139"
$ws.Range("E33").Value = "  -0.86%  "

$ws.Range("E34").Value = "  -1.19%  "

$ws.Range("D35").Value = "32.61"
$ws.Range("E35").Value = "  +10.30%  "

$ws.Range("D36").Value = "0.176"
$ws.Range("E36").Value = "  -2.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.560"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.71%  "

$ws.Range("D38").Value = "566.51"
$ws.Range("E38").Value = "  -1.07%  "

$ws.Range("D39").Value = "8.07"
$ws.Range("E39").Value = "  +7.19%  "

$ws.Range("D40").Value = "1.47"
$ws.Range("E40").Value = "  +1.99%  "

$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.150"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.15%  "

$ws.Range("D43").Value = "0.919"
$ws.Range("E43").Value = "  +0.42%  "

$ws.Range("D44").Value = "36.81"
$ws.Range("E44").Value = "  +48.08%  "

$ws.Range("D45").Value = "1.74"
$ws.Range("E45").Value = "  +2.36%  "

$ws.Range("D46").Value = "23.69"
$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("D47").Value = "5.67"
$ws.Range("E47").Value = "  +3.12%  "

$ws.Range("D48").Value = "2.27"
$ws.Range("E48").Value = "  +6.83%  "

$ws.Range("D49").Value = "0.0414"
$ws.Range("E49").Value = "  -1.66%  "

$ws.Range("D50").Value = "3.51"
$ws.Range("E50").Value = "  -1.47%  "

$ws.Range("D51").Value = "53.51"
$ws.Range("E51").Value = "  +0.59%  "
